$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "59.636.97"
Set-TextValue $ws.Range("E2") "  +0.79%  "
Set-TextValue $ws.Range("D3") "2.616.39"
Set-TextValue $ws.Range("E3") "  +0.96%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "538.27"
Set-TextValue $ws.Range("E5") "  +2.62%  "
Set-TextValue $ws.Range("D6") "142.36"
Set-TextValue $ws.Range("E6") "  +1.75%  "
Set-TextValue $ws.Range("E7") "  +0.12%  "
Set-TextValue $ws.Range("E8") "  +0.38%  "
Set-TextValue $ws.Range("E9") "  +0.53%  "
Set-TextValue $ws.Range("E10") "  +1.08%  "
Set-TextValue $ws.Range("D11") "0.335"
Set-TextValue $ws.Range("E11") "  +1.11%  "
Set-TextValue $ws.Range("E12") "  -1.33%  "
Set-TextValue $ws.Range("D13") "3.075.81"
Set-TextValue $ws.Range("E13") "  +0.73%  "
Set-TextValue $ws.Range("D14") "59.559.31"
Set-TextValue $ws.Range("E14") "  +0.79%  "
Set-TextValue $ws.Range("D15") "20.77"
Set-TextValue $ws.Range("E15") "  +1.14%  "
Set-TextValue $ws.Range("B16") "WrappedEther"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D16") "2.616.20"
Set-TextValue $ws.Range("E16") "  +1.69%  "
Set-TextValue $ws.Range("B17") "ShibaInu"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D17") "0.0000134"
Set-TextValue $ws.Range("E17") "  +0.52%  "
Set-TextValue $ws.Range("D18") "340.73"
Set-TextValue $ws.Range("E18") "  -0.40%  "
Set-TextValue $ws.Range("E19") "  +1.01%  "
Set-TextValue $ws.Range("E20") "  +0.12%  "
Set-TextValue $ws.Range("E21") "  -1.39%  "
Set-TextValue $ws.Range("E22") "  -0.07%  "
Set-TextValue $ws.Range("D23") "67.25"
Set-TextValue $ws.Range("E24") "  +0.79%  "
Set-TextValue $ws.Range("E25") "  -1.36%  "
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  +0.26%  "
Set-TextValue $ws.Range("E27") "  +2.24%  "
Set-TextValue $ws.Range("D28") "0.0₃0747"
Set-TextValue $ws.Range("E28") "  +2.86%  "
Set-TextValue $ws.Range("E30") "  +5.00%  "
Set-TextValue $ws.Range("E31") "  -1.50%  "
Set-TextValue $ws.Range("E32") "  +0.49%  "
Set-TextValue $ws.Range("D33") "150.78"
Set-TextValue $ws.Range("E33") "  +1.04%  "
Set-TextValue $ws.Range("D34") "3.99"
Set-TextValue $ws.Range("E34") "  +0.24%  "
Set-TextValue $ws.Range("E35") "  +0.37%  "
Set-TextValue $ws.Range("D36") "0.835"
Set-TextValue $ws.Range("E36") "  +2.61%  "
Set-TextValue $ws.Range("E37") "  -1.25%  "
Set-TextValue $ws.Range("E38") "  -0.26%  "
Set-TextValue $ws.Range("E39") "  +0.62%  "
Set-TextValue $ws.Range("D40") "278.60"
Set-TextValue $ws.Range("E40") "  +2.20%  "
Set-TextValue $ws.Range("D41") "0.999"
Set-TextValue $ws.Range("E41") "  +0.21%  "
Set-TextValue $ws.Range("D42") "0.602"
Set-TextValue $ws.Range("E42") "  +0.68%  "
Set-TextValue $ws.Range("D43") "10.74"
Set-TextValue $ws.Range("E43") "  -0.33%  "
Set-TextValue $ws.Range("E44") "  -0.20%  "
Set-TextValue $ws.Range("E45") "  +1.66%  "
Set-TextValue $ws.Range("D46") "1.954.02"
Set-TextValue $ws.Range("E46") "  -0.80%  "
Set-TextValue $ws.Range("E47") "  +0.33%  "
Set-TextValue $ws.Range("D48") "18.47"
Set-TextValue $ws.Range("E48") "  +1.29%  "
Set-TextValue $ws.Range("E49") "  -1.45%  "
Set-TextValue $ws.Range("D50") "111.71"
Set-TextValue $ws.Range("E50") "  -3.47%  "
Set-TextValue $ws.Range("D51") "4.74"
Set-TextValue $ws.Range("E51") "  +0.72%  "
